$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_CAN_THO")

# Update last_edited_time (column D) for rows 4,5,6,8,12,13
$rows = @(4,5,6,8,12,13)
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "2024-07-17T12:15:00.000Z"
}

# Update numeric aggregate columns on row 13 (recalculated values)
$ws.Range("T13").Value = 25800000
$ws.Range("W13").Value = 37942000
$ws.Range("AA13").Value = 112988000
$ws.Range("AE13").Value = 150930000
$ws.Range("AH13").Value = 128230000
$ws.Range("AK13").Value = 23
$ws.Range("AN13").Value = 22700000
$ws.Range("AQ13").Value = 154030000
